$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6; Excel copies formatting down from the
# row that follows (old row 6, which keeps its "HTML" category formatting),
# so A6:E6 already inherit the correct cell styles (3/3/3/3/4).
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new "corriger les liens css" entry.
$ws.Range("A6").Value = "HTML"
$ws.Range("B6").Value = "Sue la page2.html les liens css ne sont pas valides"
$ws.Range("C6").Value = "les liens pointent vers .min.css"
$ws.Range("E6").Value = "Corriger les liens"

# D6 stays empty for this entry - clear it completely so no stray cell remains.
$ws.Range("D6").Clear()

# Row 6 is a short, single-line note (unlike the 90pt row it displaced).
$ws.Rows.Item(6).RowHeight = 30

# Reflect the edit location in the view (scrolled back to the top, new cell selected).
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E6").Select()
